$d = $word.ActiveDocument

# --- Step 1: "Submitted online (tagged as 3.1)/not relevant..." -> "...3.3)/not relevant..." ---
# Surgical edit that keeps the _GoBack bookmark sitting exactly where it was (right
# after the lone "3" run), instead of doing a blind Find/Replace over the whole
# phrase (which would swallow the bookmark).
$bm = $d.Bookmarks.Item("_GoBack")
$bmPos = $bm.Start

# Remove the leading ".1" that sits right after the bookmark (turns
# ".1)/not relevant..." into ")/not relevant...").
$toDelete = $d.Range($bmPos, $bmPos + 2)
$toDelete.Text = ""

# Insert ".3" right before the bookmark (turns "...tagged as 3" into
# "...tagged as 3.3").
$insPos = $d.Range($bmPos, $bmPos)
$insPos.InsertBefore(".3")

# --- Step 2: append Questions 5, 6 and 7 after the Question 4 paragraph ---
# Paragraph 32 is "Submitted online (tagged as 3.3)/not relevant for this
# document." -- the new content goes right after it (and before the
# document's trailing empty paragraph, which must stay last).
#
# Quirks of this host worked around here:
#  - a Range with Start == End ("collapsed") sitting exactly on a paragraph
#    boundary resolves its inserts/reads ambiguously, so every "move to a
#    paragraph edge" step goes through Paragraph.Range + Collapse(...)
#    (never through $d.Range(samePos, samePos)).
#  - once a paragraph already holds text, Collapse(0) on it (its End, which
#    equals the *next* paragraph's Start) resolves to the next paragraph, so
#    follow-up appends into the same paragraph instead target End-1 (a
#    position that is unambiguously inside the paragraph's own text).

$idx = 32

function Add-Para([string]$text, [bool]$bold) {
    $tail = $d.Paragraphs.Item($idx).Range
    $tail.Collapse(0)
    $tail.InsertParagraphAfter()

    $idx = $idx + 1

    if ($bold) {
        $boldR = $d.Paragraphs.Item($idx).Range
        $boldR.Font.Bold = 1
    }

    if ($text -ne "") {
        $insertion = $d.Paragraphs.Item($idx).Range
        $insertion.Collapse(1)
        $insertion.InsertAfter($text)
    }
}

function Append-Run([string]$text) {
    # Appends more text to the paragraph most recently created by Add-Para,
    # landing it in the SAME paragraph (End-1 trick -- see note above).
    $endPos = $d.Paragraphs.Item($idx).Range.End
    $insertion = $d.Range($endPos - 1, $endPos - 1)
    $insertion.InsertAfter($text)
}

# empty, empty
Add-Para "" $false
Add-Para "" $false

# "Question 5:" (bold)
Add-Para "Question 5:" $true

# empty
Add-Para "" $false

# RMI paragraph
Add-Para "Using RMI to implement the City Server proved itself to be very straightforward and I was able to put everything up and running with minimum effort. As the RMI architecture operates on top of sockets, it is on a higher level. Instead of having to explicitly send bytes through the network I only dealt with Java objects and its methods invocations. It is important to note that some overhead is expected and sockets should be the choice for ultimate performance." $false

# empty, empty
Add-Para "" $false
Add-Para "" $false

# "Question 6:" (bold)
Add-Para "Question 6:" $true

# empty
Add-Para "" $false

# Paragraph about the DB query (3 runs of identical formatting -> inserted as one block)
Add-Para "The code was submitted online (tagged as 4.1).  " $false
Append-Run "For this task it was necessary to query the Database "
Append-Run "to define the attributes for the City object the Server was initiated with. Therefore, the client from Question 4 could be used as soon as the server is running, as the client from Question 3 is not needed anymore. The information is queried from the Database."

# empty
Add-Para "" $false

# Paragraph about JDBC changes
Add-Para "The changes necessary were basically boilerplate regarding JDBC connectivity, building an statement based on the city passed as argument to the server, and querying the DB for the information about the city in question. This way, the remote object had its attributes ready. It is important to note that it was not implemented a way to identify cities which are not in the DB." $false

# empty
Add-Para "" $false

# Paragraph about writing permissions
Add-Para "It is also worth mentioning that as we don’t have writing permissions on DB the client from Question 3 only overwrites the object’s attributes and not the information stored on the DB." $false

# empty
Add-Para "" $false

# empty bold paragraph
Add-Para "" $true

# "Question 7" (bold, single run, no colon)
Add-Para "Question 7" $true

# trailing empty bold paragraph
Add-Para "" $true

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
